# edit.ps1 — apply the CSResume.docx changes described by the diff / commit:
#   "changed edx title to intern and added programming languages to
#    spring 2016 classes"
#
# Word COM-interop (iron_native) script. $word / $d resolve to the open
# ActiveDocument.

$d = $word.ActiveDocument

# Small helper: split the run boundary that runs through absolute
# character offset $pos (relative to $base) without changing anything
# visually — toggle Bold off/on over the range [$base, $base+$pos).
# Word only splits runs across a *changed* formatting value, so we flip
# Bold to the opposite value and immediately flip it back; the net
# visual effect is nil but a run boundary is left behind at $pos.
function Split-RunBoundary($base, $pos) {
    $r = $d.Range($base, $base + $pos)
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "Probability Theory and Combinatorics, Computer Architecture
#     (languages: C, MIPS), Econometric Analysis"
#     -> "Probability and Combinatorics, Computer Architecture (C, MIPS,
#         Spark), Econometrics (Stata, GRETL, Excel)"
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("Probability Theory and Combinatorics, Computer Architecture (languages: C, MIPS), Econometric Analysis")
if ($found) {
    $start = $full.Start
    $full.Text = "Probability and Combinatorics, Computer Architecture (C, MIPS, Spark), Econometrics (Stata, GRETL, Excel)"

    # run boundaries (character offsets from $start) matching:
    #  "Probability " | "and Combinatorics, Computer Architecture" |
    #  " (" | "C, MIPS" | ", " | "Spark" | ")" |
    #  ", Econometrics (Stata, GRETL, Excel)"
    foreach ($b in @(12, 52, 54, 61, 63, 68, 69)) {
        Split-RunBoundary $start $b
    }
}

# ---------------------------------------------------------------------
# 2) "Programming languages: Python, Java, J.S., Scheme, SQL"
#     -> "Programming tools: Java, Python, JavaScript, Scheme, SQL"
# ---------------------------------------------------------------------
$full2 = $d.Content
$found2 = $full2.Find.Execute("Programming languages: Python, Java, J.S., Scheme, SQL")
if ($found2) {
    $start2 = $full2.Start
    $full2.Text = "Programming tools: Java, Python, JavaScript, Scheme, SQL"

    # run boundaries matching:
    #  "Programming tools: " | "Java," | " Python," | " " | "J" |
    #  "avaScript" | ", " | "Scheme, " | "SQL"
    foreach ($b in @(19, 24, 32, 33, 34, 43, 45, 53)) {
        Split-RunBoundary $start2 $b
    }
}

# ---------------------------------------------------------------------
# 3) "Software Developer, May 2015 – Present"
#     -> "Software Developer (Intern), May 2015 – Present"
# ---------------------------------------------------------------------
$full3 = $d.Content
$found3 = $full3.Find.Execute("Developer")
if ($found3) {
    $devStart = $full3.Start
    $devEnd = $full3.End
    $full3.Collapse(0)
    $full3.InsertAfter(" (Intern)")
    $full3.Font.Bold = 1

    # split "Developer" from the new " (Intern)" run (both bold already,
    # so force the boundary the same way as above).
    $devRange = $d.Range($devStart, $devEnd)
    $devRange.Bold = 0
    $devRange.Bold = 1
}

# ---------------------------------------------------------------------
# 4) "pedagogical planning" + "/" + "content development."
#     -> single run "pedagogical planning/content development."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("pedagogical planning/content development.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "pedagogical planning/content development.", 2) | Out-Null

Write-Output "done"
